$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 34 (shifts existing rows 34-63 down to 35-64,
# matching the original un-edited rows moving down one position).
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new Damasco/Albaricoque entry.
$ws.Range("A34").Value = 10
$ws.Range("B34").Value = "Vega Modelo de Temuco"
$ws.Range("C34").Value = "La Araucanía"
$ws.Range("D34").Value = 44893
$ws.Range("E34").Value = 9
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100103
$ws.Range("H34").Value = "Frutos de hueso (carozo)"
$ws.Range("I34").Value = 100103003
$ws.Range("J34").Value = "Damasco"
$ws.Range("K34").Value = "Albaricoque"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 250
$ws.Range("N34").Value = 15000
$ws.Range("O34").Value = 15000
$ws.Range("P34").Value = 15000
$ws.Range("Q34").Value = "$/bandeja 10 kilos"
$ws.Range("R34").Value = "Provincia de Quillota"
$ws.Range("S34").Value = 1500
$ws.Range("T34").Value = 10
